$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Simple text renames (rows above the inserted row, unaffected by the shift)
$ws.Range("A2").Value  = "Honda HR-V 2022 "
$ws.Range("A7").Value  = "Mercedes-Benz C-Class 2022 "
$ws.Range("A28").Value = "Mercedes-Benz GLC 2022 "
$ws.Range("A32").Value = "Renault Megane E-Tech 2022 "

# 2) Insert a new row at 36 (shifts existing rows 36.. down to 37..)
#    and populate it with a duplicate of the "VW Polo 2022" row (currently row 33)
$ws.Rows("36:36").Insert()
$ws.Range("A33:J33").Copy()
$ws.Range("A36").PasteSpecial()

# 3) Text renames on rows that moved down because of the inserted row
$ws.Range("A42").Value = "Isuzu D-MAX Crew Cab 2022 "
$ws.Range("A55").Value = "Mercedes-EQ EQE 2022 "
$ws.Range("A59").Value = "MAZDA CX-60 2022 "
$ws.Range("A64").Value = "Mercedes-Benz T-Class 2022 "
